$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 'https://www.ebay.com/itm/296937217728?mkevt=1&mkcid=1&mkrid=711-53200-19255-0&campid=5337590632&customid=&toolid=10050&amdata=enc%3AAQAJAAAAoHoV3kP08IDx%2BKZ9MfhVJKnEz9zHNm5Ys%2FHdWIwSYlOuE0BeX6GH42quGE3sSl5dAlam%2FqZvibRWLVcjBwKi6IUOaFV5GX3l%2Bj04s5Aot6EAd%2BGaJiJr1c08bhUqXfdxCbYr5QRV3kSai%2BKpxHvLj9JelxSO5EdWyffuJcerVyzA0Hh9pxKEADS42qEi2F7JzOTMSn1hqO39yPnOlDpVZkc%3D'
$ws.Range("N3").Value = 'https://www.ebay.com/itm/296937217728?mkevt=1&mkcid=1&mkrid=711-53200-19255-0&campid=5337590632&customid=&toolid=10050&amdata=enc%3AAQAJAAAAoHoV3kP08IDx%2BKZ9MfhVJKnEz9zHNm5Ys%2FHdWIwSYlOuE0BeX6GH42quGE3sSl5dAiCxnXDmr9IV5hd3E1C2KBzXQ8BL2C3efD%2BEOYsq4nCfYTluATd%2BT%2BiBSqCSMQaMGv4bdpkqqyx7N5JDQB96OER4dlIseJjWB3Ftt6t66VlOkbccAwuz5cnapUpiwG6j6uLxgdqf7F8VJGhNZyhN2rw%3D'
$ws.Range("C20").Value = '["https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40294/2__66737.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40293/Silver__85175.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40290/Black__60361.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40295/Green__35373.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40287/4__70292.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40288/3__05149.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40289/1__14954.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/224474/newmainpyp_54_62aa0b972b812987f4bda9f9__90767.1738031866.png?c=1"]'
$ws.Range("J20").Value = '["https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40294/2__66737.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40293/Silver__85175.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40290/Black__60361.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40295/Green__35373.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40287/4__70292.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40288/3__05149.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/40289/1__14954.1738031866.jpg?c=1", "https://cdn11.bigcommerce.com/s-yqzpac968l/images/stencil/1280x1280/products/1954/224474/newmainpyp_54_62aa0b972b812987f4bda9f9__90767.1738031866.png?c=1"]'
